# ---------------------------------------------------------------------------
# Applies the "refinements to metrics and plot outputs + new month run" edit:
#   1. Country sheet: update CONCERN.HIGH / MoM figures for the latest run.
#   2. States sheet: insert a STATE.CODE column (ISO codes) in front of
#      STATE.NAME, re-derive CONCERN.HIGH ordering for the new month, and
#      refresh MoM deltas.
#   3. Dark clusters sheet: update the Northeast cluster figures.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Country sheet
# ---------------------------------------------------------------------------
$country = $wb.Worksheets.Item("Country")
$country.Range("A2").Value = 32.78
$country.Range("B2").Value = 3.017

# ---------------------------------------------------------------------------
# 2) States sheet
# ---------------------------------------------------------------------------
$states = $wb.Worksheets.Item("States")

# Insert a new first column for the ISO state code; existing STATE.NAME /
# CONCERN.HIGH / MoM columns all shift one place to the right.
$states.Columns.Item(1).Insert()

$states.Range("A1").Value = "STATE.CODE"
$states.Range("B1").Value = "STATE.NAME"
$states.Range("C1").Value = "CONCERN.HIGH"
$states.Range("D1").Value = "MoM"

$statesData = New-Object 'object[,]' 32,4

function Set-StateRow($idx, $code, $name, $concern, $mom) {
    $statesData[$idx,0] = $code
    $statesData[$idx,1] = $name
    $statesData[$idx,2] = $concern
    $statesData[$idx,3] = $mom
}

Set-StateRow 0  "IN-BR" "Bihar"                    89.47 0
Set-StateRow 1  "IN-MZ" "Mizoram"                  81.82 49.9908
Set-StateRow 2  "IN-JH" "Jharkhand"                75    0
Set-StateRow 3  "IN-PB" "Punjab"                   63.64 0
Set-StateRow 4  "IN-ML" "Meghalaya"                63.64 0
Set-StateRow 5  "IN-UP" "Uttar Pradesh"            62.67 0
Set-StateRow 6  "IN-MN" "Manipur"                  62.5  11.1111
Set-StateRow 7  "IN-NL" "Nagaland"                 54.55 50.0275
Set-StateRow 8  "IN-TR" "Tripura"                  50    0
Set-StateRow 9  "IN-DD" "Daman and Diu"            50    0
Set-StateRow 10 "IN-PY" "Puducherry"               50    0
Set-StateRow 11 "IN-DL" "Delhi"                    45.45 0
Set-StateRow 12 "IN-AR" "Arunachal Pradesh"        44    0
Set-StateRow 13 "IN-OR" "Odisha"                   33.33 11.1
Set-StateRow 14 "IN-TS" "Telangana"                33.33 0
Set-StateRow 15 "IN-HR" "Haryana"                  31.82 0
Set-StateRow 16 "IN-MP" "Madhya Pradesh"           30.77 0
Set-StateRow 17 "IN-JK" "Jammu and Kashmir"        27.27 0
Set-StateRow 18 "IN-AS" "Assam"                    24.24 0
Set-StateRow 19 "IN-HP" "Himachal Pradesh"         16.67 0
Set-StateRow 20 "IN-RJ" "Rajasthan"                15.15 0
Set-StateRow 21 "IN-WB" "West Bengal"              13.04 0
Set-StateRow 22 "IN-CT" "Chhattisgarh"             7.41  0
Set-StateRow 23 "IN-GJ" "Gujarat"                  3.03  0
Set-StateRow 24 "IN-MH" "Maharashtra"              2.78  0
Set-StateRow 25 "IN-CH" "Chandigarh"               0     0
Set-StateRow 26 "IN-DN" "Dadra and Nagar Haveli"   0     0
Set-StateRow 27 "IN-LA" "Ladakh"                   0     0
Set-StateRow 28 "IN-AP" "Andhra Pradesh"           0     0
Set-StateRow 29 "IN-TN" "Tamil Nadu"               0     0
Set-StateRow 30 "IN-KA" "Karnataka"                0     0
Set-StateRow 31 "IN-UL" "Uttarakhand"              0     0

$states.Range("A2:D33").Value = $statesData

# The last seven states (Chandigarh .. Uttarakhand) never had a MoM value in
# the original sheet - keep column D blank for them, matching the source.
$states.Range("D27:D33").ClearContents()

# ---------------------------------------------------------------------------
# 3) Dark clusters sheet
# ---------------------------------------------------------------------------
$dark = $wb.Worksheets.Item("Dark clusters")
$dark.Range("B3").Value = 47.83
$dark.Range("C3").Value = 12.2506
